$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 - title slide: give the ctrTitle and subTitle placeholders
# explicit positions/sizes (previously inherited from the layout).
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$s1Title = $s1.Shapes.Item(1)
$s1Title.Left = 150.79748031496064
$s1Title.Top = 176.473937007874
$s1Title.Width = 658.3644881889763
$s1Title.Height = 165.21465366929135

$s1Sub = $s1.Shapes.Item(2)
$s1Sub.Left = 211.01622847244096
$s1Sub.Top = 349.2655905511811
$s1Sub.Width = 537.9270338740158
$s1Sub.Height = 85.53047644094488

# ---------------------------------------------------------------------
# Slide 3 - "Purpose": bump the line spacing of all three bullet
# paragraphs to 150%.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3Body = $s3.Shapes.Item(2)
$s3Paras = $s3Body.TextFrame.TextRange.Paragraphs()
$s3Paras.ParagraphFormat.SpaceWithin = 1.5

# ---------------------------------------------------------------------
# Slide 4 - "Data Science process" -> "Data Science Process"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Data Science Process"

# ---------------------------------------------------------------------
# Slide 8 - "Comparison of the year on year closing price": give the
# title an explicit position/size and reposition/resize the picture.
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)

$s8Title = $s8.Shapes.Item(1)
$s8Title.Left = 95.41740157480315
$s8Title.Top = 28.83496062992126
$s8Title.Width = 841.9805611811023
$s8Title.Height = 117.0

$s8Pic = $s8.Shapes.Item(2)
$s8Pic.Left = 302.05110236220474
$s8Pic.Top = 137.07299212598426
$s8Pic.Width = 394.8809448818898
$s8Pic.Height = 343.10182302362205

# ---------------------------------------------------------------------
# Slide 9 - "The correlation between Open and High" ->
# "The Correlation between Open and High"; reposition/resize picture.
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = "The Correlation between Open and High"

$s9Pic = $s9.Shapes.Item(2)
$s9Pic.Left = 285.7411971023622
$s9Pic.Top = 139.80582677165356
$s9Pic.Width = 388.5175630551181
$s9Pic.Height = 337.5728346456693

# ---------------------------------------------------------------------
# Slide 10 - "Linear regression" -> "Linear Regression"
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Linear Regression"

# ---------------------------------------------------------------------
# Slide 11 - "Prediction result" -> "Prediction Result"
# ---------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = "Prediction Result"

# ---------------------------------------------------------------------
# Slide 12 - "ARIMA" -> "Forecasting with ARIMA"
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = "Forecasting with ARIMA"
